$d = $word.ActiveDocument

# Insert a new centered paragraph containing the contact info directly
# after the "Dheeraj Chand" name line, using Find/Replace with a
# paragraph-mark (^p) in the replacement text. This creates a clean new
# paragraph that only inherits the paragraph-level centering from the
# matched paragraph, without dragging along the bold/28pt run formatting
# that "Dheeraj Chand" uses.
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
